$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# "animate enemy movement" and "enemy color?" were finished and replaced by
# a single "enemy sprites" to-do item; the list shifts up and two new
# pause-menu related to-do items are appended, while the old last row
# ("track score with board") is removed entirely.
$ws.Range("A3").Value = "enemy sprites"
$ws.Range("A4").Value = "change direction"
$ws.Range("A5").Value = "game name on pause menu"
$ws.Range("A6").Value = "align lesser congratulatory message"
$ws.Range("A7").ClearContents()

$ws.Range("B8").Select() | Out-Null
